$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block: rows 16-19 (n-splitting check) ---
$ws.Range("D16").Value = -0.742
$ws.Range("E16").Value = 0.756
$ws.Range("G16").Value = 14.8
$ws.Range("H16").Value = 14.2

$ws.Range("D17").Value = 0.227
$ws.Range("E17").Value = 0.475
$ws.Range("G17").Value = 1.42
$ws.Range("H17").Value = 1.01

$ws.Range("D18").Formula = "=ABS(D16-D17)"
$ws.Range("E18").Formula = "=SQRT(E16^2+E17^2)"
$ws.Range("G18").Formula = "=ABS(G16-G17)"
$ws.Range("H18").Formula = "=SQRT(H16^2+H17^2)"

$ws.Range("D19").Formula = "=D18/E18"
$ws.Range("G19").Formula = "=G18/H18"

# --- New data block: rows 29-35 (stdev check) ---
$ws.Range("D29").Value = 0.77299340000000005
$ws.Range("D30").Value = -0.98535720000000004
$ws.Range("D31").Value = 0.23331389999999999
$ws.Range("D32").Value = 0.37578810000000001
$ws.Range("D33").Value = -0.7940083
$ws.Range("D34").Value = 1.7285349999999999

$ws.Range("D35").Formula = "=STDEV(D29:D34)"

# --- Update selection / active cell ---
$ws.Range("D35").Select()

# --- Move/resize the chart graphic frame ---
# Target anchor (0-based col/row + EMU offset, matching the OOXML <xdr:from>/<xdr:to>):
#   from: col=10 colOff=488950 row=21 rowOff=88900
#   to:   col=16 colOff=654050 row=34 rowOff=190500
# ChartObjects(i).Left/Top/Width/Height are in points and drive the anchor,
# so convert the target cell-anchor (0-based col/row + EMU offset) to points
# using the sheet's actual column/row geometry.
function Get-AnchorPoint($sheet, $col0, $colOff, $row0, $rowOff) {
    $ptLeft = $sheet.Cells.Item(1, $col0 + 1).Left + ($colOff / 12700.0)
    $ptTop = $sheet.Cells.Item($row0 + 1, 1).Top + ($rowOff / 12700.0)
    return @{ Left = $ptLeft; Top = $ptTop }
}

$fromPt = Get-AnchorPoint $ws 10 488950 21 88900
$toPt = Get-AnchorPoint $ws 16 654050 34 190500

$co = $ws.ChartObjects(1)
$co.Left = $fromPt.Left
$co.Top = $fromPt.Top
$co.Width = $toPt.Left - $fromPt.Left
$co.Height = $toPt.Top - $fromPt.Top
